$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '中国长城'
$ws.Range("B2").Value = '利欧股份'
$ws.Range("C2").Value = '国晟科技'
$ws.Range("A3").Value = '航天电子'
$ws.Range("B3").Value = '国晟科技'
$ws.Range("C3").Value = '利欧股份'
$ws.Range("A4").Value = '通富微电'
$ws.Range("B4").Value = '白银有色'
$ws.Range("C4").Value = '中国长城'
$ws.Range("A5").Value = '利欧股份'
$ws.Range("B5").Value = '湖南白银'
$ws.Range("C5").Value = '锋龙股份'
$ws.Range("A6").Value = '国晟科技'
$ws.Range("B6").Value = '中国长城'
$ws.Range("C6").Value = '湖南白银'
$ws.Range("A7").Value = '巨力索具'
$ws.Range("B7").Value = '航天电子'
$ws.Range("C7").Value = '航天电子'
$ws.Range("A8").Value = '锋龙股份'
$ws.Range("B8").Value = '锋龙股份'
$ws.Range("C8").Value = '白银有色'
$ws.Range("A9").Value = '白银有色'
$ws.Range("B9").Value = '通富微电'
$ws.Range("C9").Value = '雪人集团'
$ws.Range("A10").Value = '湖南白银'
$ws.Range("B10").Value = '华天科技'
$ws.Range("C10").Value = '巨力索具'
$ws.Range("A11").Value = '信维通信'
$ws.Range("B11").Value = '金风科技'
$ws.Range("C11").Value = '雷科防务'
$ws.Range("A12").Value = '雪人集团'
$ws.Range("B12").Value = '雪人集团'
$ws.Range("C12").Value = '金风科技'
$ws.Range("A13").Value = '雷科防务'
$ws.Range("B13").Value = '雷科防务'
$ws.Range("C13").Value = '航天发展'
$ws.Range("A14").Value = '金风科技'
$ws.Range("B14").Value = '岩山科技'
$ws.Range("C14").Value = '通富微电'
$ws.Range("A15").Value = '航天发展'
$ws.Range("B15").Value = '巨力索具'
$ws.Range("C15").Value = '中国卫通'
$ws.Range("A16").Value = '岩山科技'
$ws.Range("B16").Value = '航天发展'
$ws.Range("C16").Value = '嘉美包装'
$ws.Range("A17").Value = '中国卫通'
$ws.Range("B17").Value = '盈方微'
$ws.Range("C17").Value = '岩山科技'
$ws.Range("A18").Value = '华天科技'
$ws.Range("B18").Value = '中国西电'
$ws.Range("C18").Value = '中超控股'
$ws.Range("A19").Value = '海光信息'
$ws.Range("B19").Value = '三角防务'
$ws.Range("C19").Value = '海格通信'
$ws.Range("A20").Value = '金安国纪'
$ws.Range("B20").Value = '中超控股'
$ws.Range("C20").Value = '久其软件'
$ws.Range("A21").Value = '西部材料'
$ws.Range("B21").Value = '中国卫通'
$ws.Range("C21").Value = '长电科技'
